$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 7. This pushes the previous rows 7..113
# down to 8..114, so the whole table grows by one (weekly) record.
$ws.Rows(7).Insert()

# Copy across the static/unchanged columns from the row that just got shifted
# down into row 8 (this used to be row 7 before the insert), then overwrite
# the columns that hold this week's new figures.
$ws.Range("A8:C8").Copy()
$ws.Range("A7").PasteSpecial(-4104)  # xlPasteAll
$ws.Range("E8:I8").Copy()
$ws.Range("E7").PasteSpecial(-4104)
$ws.Range("N8:R8").Copy()
$ws.Range("N7").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Cells.Item(7, 4).Value = 44750   # D7 Fecha
$ws.Cells.Item(7, 9).Value = "Primera"  # I7 Calidad
$ws.Cells.Item(7, 10).Value = 140    # J7 Volumen
$ws.Cells.Item(7, 11).Value = 33000  # K7 Precio minimo
$ws.Cells.Item(7, 12).Value = 35000  # L7 Precio maximo
$ws.Cells.Item(7, 13).Value = 34000  # M7 Precio promedio ponderado
$ws.Cells.Item(7, 16).Value = 1700   # P7 Precio $/Kg
